$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Sema3e"
$ws.Cells.Item(2, 3).Value = "Plxnd1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.1244045
$ws.Cells.Item(2, 8).Value = 0.248809
$ws.Cells.Item(2, 9).Value = 0.02044382138936323
$ws.Cells.Item(2, 10).Value = 0.0203137528521298
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 95.41999849999999
$ws.Cells.Item(2, 14).Value = 190.839997
$ws.Cells.Item(2, 15).Value = 0.4569018340225213
$ws.Cells.Item(2, 16).Value = 0.3855869467978643
$ws.Cells.Item(2, 17).Value = 11.87067720339325
$ws.Cells.Item(2, 18).Value = 47.482708813573
$ws.Cells.Item(2, 19).Value = 0.00934081948722891
$ws.Cells.Item(2, 20).Value = 0.007832717940259139

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Sema3e"
$ws.Cells.Item(3, 3).Value = "Plxnd1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.1244045
$ws.Cells.Item(3, 8).Value = 0.248809
$ws.Cells.Item(3, 9).Value = 0.02044382138936323
$ws.Cells.Item(3, 10).Value = 0.0203137528521298
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 16.807086
$ws.Cells.Item(3, 14).Value = 50.42125799999999
$ws.Cells.Item(3, 15).Value = 0.08047776712105316
$ws.Cells.Item(3, 16).Value = 0.1018747601737145
$ws.Cells.Item(3, 17).Value = 2.090877130287
$ws.Cells.Item(3, 18).Value = 12.545262781722
$ws.Cells.Item(3, 19).Value = 0.001645273096837579
$ws.Cells.Item(3, 20).Value = 0.002069458700038832

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Sema3e"
$ws.Cells.Item(4, 3).Value = "Plxnd1"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.1244045
$ws.Cells.Item(4, 8).Value = 0.248809
$ws.Cells.Item(4, 9).Value = 0.02044382138936323
$ws.Cells.Item(4, 10).Value = 0.0203137528521298
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 24.47797766666666
$ws.Cells.Item(4, 14).Value = 73.433933
$ws.Cells.Item(4, 15).Value = 0.117208478986324
$ws.Cells.Item(4, 16).Value = 0.1483712348666036
$ws.Cells.Item(4, 17).Value = 3.045170572632833
$ws.Cells.Item(4, 18).Value = 18.271023435797
$ws.Cells.Item(4, 19).Value = 0.002396189209715341
$ws.Cells.Item(4, 20).Value = 0.00301397659544549

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Sema3e"
$ws.Cells.Item(5, 3).Value = "Plxnd1"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.1244045
$ws.Cells.Item(5, 8).Value = 0.248809
$ws.Cells.Item(5, 9).Value = 0.02044382138936323
$ws.Cells.Item(5, 10).Value = 0.0203137528521298
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 36.17031
$ws.Cells.Item(5, 14).Value = 72.34062
$ws.Cells.Item(5, 15).Value = 0.1731951502405771
$ws.Cells.Item(5, 16).Value = 0.1461622261252945
$ws.Cells.Item(5, 17).Value = 4.499749330395
$ws.Cells.Item(5, 18).Value = 17.99899732158
$ws.Cells.Item(5, 19).Value = 0.003540770717022289
$ws.Cells.Item(5, 20).Value = 0.002969103337826342

$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Sema3e"
$ws.Cells.Item(6, 3).Value = "Plxnd1"
$ws.Cells.Item(6, 4).Value = "Neutrophils"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.1244045
$ws.Cells.Item(6, 8).Value = 0.248809
$ws.Cells.Item(6, 9).Value = 0.02044382138936323
$ws.Cells.Item(6, 10).Value = 0.0203137528521298
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 8.119321333333334
$ws.Cells.Item(6, 14).Value = 24.357964
$ws.Cells.Item(6, 15).Value = 0.03887793823658658
$ws.Cells.Item(6, 16).Value = 0.04921459398771786
$ws.Cells.Item(6, 17).Value = 1.010080110812667
$ws.Cells.Item(6, 18).Value = 6.060480664876
$ws.Cells.Item(6, 19).Value = 0.0007948136252954712
$ws.Cells.Item(6, 20).Value = 0.0009997330989844138

$ws.Cells.Item(7, 1).Value = "ECs"
$ws.Cells.Item(7, 2).Value = "Sema3e"
$ws.Cells.Item(7, 3).Value = "Plxnd1"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.1244045
$ws.Cells.Item(7, 8).Value = 0.248809
$ws.Cells.Item(7, 9).Value = 0.02044382138936323
$ws.Cells.Item(7, 10).Value = 0.0203137528521298
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 27.84666233333333
$ws.Cells.Item(7, 14).Value = 83.539987
$ws.Cells.Item(7, 15).Value = 0.1333388313929376
$ws.Cells.Item(7, 16).Value = 0.1687902380488052
$ws.Cells.Item(7, 17).Value = 3.464250104247166
$ws.Cells.Item(7, 18).Value = 20.785500625483
$ws.Cells.Item(7, 19).Value = 0.002725955253263636
$ws.Cells.Item(7, 20).Value = 0.003428763179575586

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Sema3e"
$ws.Cells.Item(8, 3).Value = "Plxnd1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.07792666666666666
$ws.Cells.Item(8, 8).Value = 0.23378
$ws.Cells.Item(8, 9).Value = 0.01280595842434782
$ws.Cells.Item(8, 10).Value = 0.01908672572845397
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 95.41999849999999
$ws.Cells.Item(8, 14).Value = 190.839997
$ws.Cells.Item(8, 15).Value = 0.4569018340225213
$ws.Cells.Item(8, 16).Value = 0.3855869467978643
$ws.Cells.Item(8, 17).Value = 7.435762416443332
$ws.Cells.Item(8, 18).Value = 44.61457449865999
$ws.Cells.Item(8, 19).Value = 0.005851065890500677
$ws.Cells.Item(8, 20).Value = 0.00735959229800281

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Sema3e"
$ws.Cells.Item(9, 3).Value = "Plxnd1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.07792666666666666
$ws.Cells.Item(9, 8).Value = 0.23378
$ws.Cells.Item(9, 9).Value = 0.01280595842434782
$ws.Cells.Item(9, 10).Value = 0.01908672572845397
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 16.807086
$ws.Cells.Item(9, 14).Value = 50.42125799999999
$ws.Cells.Item(9, 15).Value = 0.08047776712105316
$ws.Cells.Item(9, 16).Value = 0.1018747601737145
$ws.Cells.Item(9, 17).Value = 1.30972018836
$ws.Cells.Item(9, 18).Value = 11.78748169524
$ws.Cells.Item(9, 19).Value = 0.001030594939836553
$ws.Cells.Item(9, 20).Value = 0.001944455606087714

$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Sema3e"
$ws.Cells.Item(10, 3).Value = "Plxnd1"
$ws.Cells.Item(10, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.07792666666666666
$ws.Cells.Item(10, 8).Value = 0.23378
$ws.Cells.Item(10, 9).Value = 0.01280595842434782
$ws.Cells.Item(10, 10).Value = 0.01908672572845397
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 24.47797766666666
$ws.Cells.Item(10, 14).Value = 73.433933
$ws.Cells.Item(10, 15).Value = 0.117208478986324
$ws.Cells.Item(10, 16).Value = 0.1483712348666036
$ws.Cells.Item(10, 17).Value = 1.907487206304444
$ws.Cells.Item(10, 18).Value = 17.16738485674
$ws.Cells.Item(10, 19).Value = 0.00150096690887991
$ws.Cells.Item(10, 20).Value = 0.002831921065890891

$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Sema3e"
$ws.Cells.Item(11, 3).Value = "Plxnd1"
$ws.Cells.Item(11, 4).Value = "MuSCs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.07792666666666666
$ws.Cells.Item(11, 8).Value = 0.23378
$ws.Cells.Item(11, 9).Value = 0.01280595842434782
$ws.Cells.Item(11, 10).Value = 0.01908672572845397
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 36.17031
$ws.Cells.Item(11, 14).Value = 72.34062
$ws.Cells.Item(11, 15).Value = 0.1731951502405771
$ws.Cells.Item(11, 16).Value = 0.1461622261252945
$ws.Cells.Item(11, 17).Value = 2.8186316906
$ws.Cells.Item(11, 18).Value = 16.9117901436
$ws.Cells.Item(11, 19).Value = 0.002217929893279505
$ws.Cells.Item(11, 20).Value = 0.002789758321913766

$ws.Cells.Item(12, 1).Value = "FAPs"
$ws.Cells.Item(12, 2).Value = "Sema3e"
$ws.Cells.Item(12, 3).Value = "Plxnd1"
$ws.Cells.Item(12, 4).Value = "Neutrophils"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.07792666666666666
$ws.Cells.Item(12, 8).Value = 0.23378
$ws.Cells.Item(12, 9).Value = 0.01280595842434782
$ws.Cells.Item(12, 10).Value = 0.01908672572845397
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 8.119321333333334
$ws.Cells.Item(12, 14).Value = 24.357964
$ws.Cells.Item(12, 15).Value = 0.03887793823658658
$ws.Cells.Item(12, 16).Value = 0.04921459398771786
$ws.Cells.Item(12, 17).Value = 0.6327116471022222
$ws.Cells.Item(12, 18).Value = 5.694404823919999
$ws.Cells.Item(12, 19).Value = 0.0004978692606820901
$ws.Cells.Item(12, 20).Value = 0.0009393454572807907

$ws.Cells.Item(13, 1).Value = "FAPs"
$ws.Cells.Item(13, 2).Value = "Sema3e"
$ws.Cells.Item(13, 3).Value = "Plxnd1"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.07792666666666666
$ws.Cells.Item(13, 8).Value = 0.23378
$ws.Cells.Item(13, 9).Value = 0.01280595842434782
$ws.Cells.Item(13, 10).Value = 0.01908672572845397
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 27.84666233333333
$ws.Cells.Item(13, 14).Value = 83.539987
$ws.Cells.Item(13, 15).Value = 0.1333388313929376
$ws.Cells.Item(13, 16).Value = 0.1687902380488052
$ws.Cells.Item(13, 17).Value = 2.169997573428888
$ws.Cells.Item(13, 18).Value = 19.52997816086
$ws.Cells.Item(13, 19).Value = 0.001707531531169083
$ws.Cells.Item(13, 20).Value = 0.003221652979278002

$ws.Cells.Item(14, 1).Value = "MuSCs"
$ws.Cells.Item(14, 2).Value = "Sema3e"
$ws.Cells.Item(14, 3).Value = "Plxnd1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 5.882857
$ws.Cells.Item(14, 8).Value = 11.765714
$ws.Cells.Item(14, 9).Value = 0.966750220186289
$ws.Cells.Item(14, 10).Value = 0.9605995214194162
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 95.41999849999999
$ws.Cells.Item(14, 14).Value = 190.839997
$ws.Cells.Item(14, 15).Value = 0.4569018340225213
$ws.Cells.Item(14, 16).Value = 0.3855869467978643
$ws.Cells.Item(14, 17).Value = 561.3422061157144
$ws.Cells.Item(14, 18).Value = 2245.368824462857
$ws.Cells.Item(14, 19).Value = 0.4417099486447917
$ws.Cells.Item(14, 20).Value = 0.3703946365596024

$ws.Cells.Item(15, 1).Value = "MuSCs"
$ws.Cells.Item(15, 2).Value = "Sema3e"
$ws.Cells.Item(15, 3).Value = "Plxnd1"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 5.882857
$ws.Cells.Item(15, 8).Value = 11.765714
$ws.Cells.Item(15, 9).Value = 0.966750220186289
$ws.Cells.Item(15, 10).Value = 0.9605995214194162
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 16.807086
$ws.Cells.Item(15, 14).Value = 50.42125799999999
$ws.Cells.Item(15, 15).Value = 0.08047776712105316
$ws.Cells.Item(15, 16).Value = 0.1018747601737145
$ws.Cells.Item(15, 17).Value = 98.87368352470199
$ws.Cells.Item(15, 18).Value = 593.2421011482119
$ws.Cells.Item(15, 19).Value = 0.07780189908437903
$ws.Cells.Item(15, 20).Value = 0.09786084586758793

$ws.Cells.Item(16, 1).Value = "MuSCs"
$ws.Cells.Item(16, 2).Value = "Sema3e"
$ws.Cells.Item(16, 3).Value = "Plxnd1"
$ws.Cells.Item(16, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 5.882857
$ws.Cells.Item(16, 8).Value = 11.765714
$ws.Cells.Item(16, 9).Value = 0.966750220186289
$ws.Cells.Item(16, 10).Value = 0.9605995214194162
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 24.47797766666666
$ws.Cells.Item(16, 14).Value = 73.433933
$ws.Cells.Item(16, 15).Value = 0.117208478986324
$ws.Cells.Item(16, 16).Value = 0.1483712348666036
$ws.Cells.Item(16, 17).Value = 144.0004422621936
$ws.Cells.Item(16, 18).Value = 864.0026535731619
$ws.Cells.Item(16, 19).Value = 0.1133113228677288
$ws.Cells.Item(16, 20).Value = 0.1425253372052672

$ws.Cells.Item(17, 1).Value = "MuSCs"
$ws.Cells.Item(17, 2).Value = "Sema3e"
$ws.Cells.Item(17, 3).Value = "Plxnd1"
$ws.Cells.Item(17, 4).Value = "MuSCs"
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 5.882857
$ws.Cells.Item(17, 8).Value = 11.765714
$ws.Cells.Item(17, 9).Value = 0.966750220186289
$ws.Cells.Item(17, 10).Value = 0.9605995214194162
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 36.17031
$ws.Cells.Item(17, 14).Value = 72.34062
$ws.Cells.Item(17, 15).Value = 0.1731951502405771
$ws.Cells.Item(17, 16).Value = 0.1461622261252945
$ws.Cells.Item(17, 17).Value = 212.78476137567
$ws.Cells.Item(17, 18).Value = 851.13904550268
$ws.Cells.Item(17, 19).Value = 0.1674364496302753
$ws.Cells.Item(17, 20).Value = 0.1404033644655544

$ws.Cells.Item(18, 1).Value = "MuSCs"
$ws.Cells.Item(18, 2).Value = "Sema3e"
$ws.Cells.Item(18, 3).Value = "Plxnd1"
$ws.Cells.Item(18, 4).Value = "Neutrophils"
$ws.Cells.Item(18, 5).Value = 2
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 5.882857
$ws.Cells.Item(18, 8).Value = 11.765714
$ws.Cells.Item(18, 9).Value = 0.966750220186289
$ws.Cells.Item(18, 10).Value = 0.9605995214194162
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 8.119321333333334
$ws.Cells.Item(18, 14).Value = 24.357964
$ws.Cells.Item(18, 15).Value = 0.03887793823658658
$ws.Cells.Item(18, 16).Value = 0.04921459398771786
$ws.Cells.Item(18, 17).Value = 47.76480634104933
$ws.Cells.Item(18, 18).Value = 286.588838046296
$ws.Cells.Item(18, 19).Value = 0.03758525535060902
$ws.Cells.Item(18, 20).Value = 0.04727551543145266

$ws.Cells.Item(19, 1).Value = "MuSCs"
$ws.Cells.Item(19, 2).Value = "Sema3e"
$ws.Cells.Item(19, 3).Value = "Plxnd1"
$ws.Cells.Item(19, 4).Value = "Resolving-Mac"
$ws.Cells.Item(19, 5).Value = 2
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 5.882857
$ws.Cells.Item(19, 8).Value = 11.765714
$ws.Cells.Item(19, 9).Value = 0.966750220186289
$ws.Cells.Item(19, 10).Value = 0.9605995214194162
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 27.84666233333333
$ws.Cells.Item(19, 14).Value = 83.539987
$ws.Cells.Item(19, 15).Value = 0.1333388313929376
$ws.Cells.Item(19, 16).Value = 0.1687902380488052
$ws.Cells.Item(19, 17).Value = 163.8179324342863
$ws.Cells.Item(19, 18).Value = 982.9075946057179
$ws.Cells.Item(19, 19).Value = 0.1289053446085049
$ws.Cells.Item(19, 20).Value = 0.1621398218899517
